$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238853812217712
$ws.Range("B1").Value = 2.996760845184326
$ws.Range("C1").Value = 6.001397132873535
$ws.Range("D1").Value = 1.950359582901001
$ws.Range("E1").Value = 0.7457773089408875
